# daily auto push: 2026-02-10 14:24 UTC
# Insert a new log row right before the "2026/12/29" block (old row 779),
# shifting all subsequent rows down by one, and populate the new row with
# today's entry (date 2026/02/10, weekday 火, hour 21, ranking 132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 779..820 down to 780..821 by inserting a blank row at 779.
$ws.Rows("779").Insert()

# Write the new row's values. The date column stores a plain text value
# (matching the rest of the sheet, which never uses real Excel date
# serials) rather than letting Excel auto-parse "2026/02/10" into a date
# serial number. A leading apostrophe forces text entry, and resetting
# the style back to "Normal" afterwards clears the quote-prefix styling
# that Excel would otherwise leave behind, keeping the cell unstyled just
# like its neighbours.
$ws.Range("A779").Value = "'2026/02/10"
$ws.Range("A779").Style = "Normal"
$ws.Range("B779").Value = "火"
$ws.Range("C779").Value = 21
$ws.Range("D779").Value = 132
